# Update the "DO NOT TOUCH AUTOMATION HMRC RECOGNITION EMPLOYEE 106" marker text
# to "DO NOT TOUCH AUTOMATION HMRC RECOGNITION EMPLOYEE 6" everywhere it appears.
$wb = $excel.ActiveWorkbook

$newText = "DO NOT TOUCH AUTOMATION HMRC RECOGNITION EMPLOYEE 6"

$wsReset = $wb.Worksheets.Item("ResetEmployeeData12")
$wsReset.Range("A2").Value = $newText

$wsTax = $wb.Worksheets.Item("UpdteTaxCodeAndAnualSalaryM12")
$wsTax.Range("A2").Value = $newText

$wsMarch = $wb.Worksheets.Item("ProcessPayrollForMarch")
$wsMarch.Range("B2").Value = $newText

$wsFinal = $wb.Worksheets.Item("ProcessFinalPayrollForMarch")
$wsFinal.Range("B2").Value = $newText

$wsTest = $wb.Worksheets.Item("TestMarchReports")
$wsTest.Range("B2").Value = $newText

# Update the remembered selection (active cell) on each sheet.
$wsFirst = $wb.Worksheets.Item("first")
$wsFirst.Range("A3").Select() | Out-Null

$wsReset.Range("A2").Select() | Out-Null
$wsTax.Range("A2").Select() | Out-Null
$wsMarch.Range("B2").Select() | Out-Null
$wsFinal.Range("B2").Select() | Out-Null

# TestMarchReports becomes the active (selected) sheet, with B2 selected.
$wsTest.Activate() | Out-Null
$wsTest.Range("B2").Select() | Out-Null

# (Tab-strip scroll position ("first visible tab") is a window-chrome detail
# with no reachable Excel object-model property in this host; activating the
# sheet above is what drives bookViews/workbookView@activeTab.)
